{"js": "const ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">High school was a lonely place for me. I bought into the competitive environment where if you studied hard and long enough and did enough extracurriculars, you would get into a good college and your future would be secure. I didn\u2019t make time for friends, and my friends didn\u2019t make time for me. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">I had a lot of friends of convenience, like track friends and AP calc friends and band friends, but very few people were simply friends. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">This brought a deep loneliness into me, and with a deep loneliness came the belief that my relationships would fix all the unhappiness that I had. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">With the friends I made through my life, I became as loving as I could because spending time with them cured the small loneliness that I had. For the short moments that I was with my friends, I felt whole, and when I was alone in my basement with nothing to do besides search the internet, I was empty. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Did I love my friends? Yes</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>,</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> I did, but only because they made me feel not empty. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">My greatest opiate was my best friend Ivy Chen. Throughout high school, she was the one constant friend in my life. I loved her for it, but my love had limits. In the summer before college, she became so busy with her job that we couldn\u2019t talk any more. I became so embittered with her because she wouldn\u2019t make space in her life for me, that we didn\u2019t talk for almost a year. </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>Coming to college, I was hoping for a clean slate and for new friends, and God provided. A random guy i</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">n my lecture invited me to JCA, and I found the friends I </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">was so desperately looking for. Not only that, but I found a community and a God that loved in a way that was so different from how I loved. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>At the same time, a lot of friends from home confided in me about brokenness in their lives, whether it was through a bad relationship or a family sickness or a simple lack of direction. For these people whom I loved, I prayed deeply. I wanted to go to morning prayer just to pray for my friends. I li</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>stened to a lot of sermons too. I thought I was doing all that I could to love and be loved.</w:t>\n      </w:r>\n      <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n      <w:bookmarkEnd w:id=\"0\"/>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> </w:t>\n      </w:r>\n    </w:p>\n</w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ncontext.document.body.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">High school was a lonely place for me. I bought into the competitive environment where if you studied hard and long enough and did enough extracurriculars, you would get into a good college and your future would be secure. I didn\u2019t make time for friends, and my friends didn\u2019t make time for me. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">I had a lot of friends of convenience, like track friends and AP calc friends and band friends, but very few people were simply friends. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">This brought a deep loneliness into me, and with a deep loneliness came the belief that my relationships would fix all the unhappiness that I had. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">With the friends I made through my life, I became as loving as I could because spending time with them cured the small loneliness that I had. For the short moments that I was with my friends, I felt whole, and when I was alone in my basement with nothing to do besides search the internet, I was empty. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>Did I love my friends? Yes</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>,</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> I did, but only because they made me feel not empty. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">My greatest opiate was my best friend Ivy Chen. Throughout high school, she was the one constant friend in my life. I loved her for it, but my love had limits. In the summer before college, she became so busy with her job that we couldn\u2019t talk any more. I became so embittered with her because she wouldn\u2019t make space in her life for me, that we didn\u2019t talk for almost a year. </w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>Coming to college, I was hoping for a clean slate and for new friends, and God provided. A random guy i</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">n my lecture invited me to JCA, and I found the friends I </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">was so desperately looking for. Not only that, but I found a community and a God that loved in a way that was so different from how I loved. </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>At the same time, a lot of friends from home confided in me about brokenness in their lives, whether it was through a bad relationship or a family sickness or a simple lack of direction. For these people whom I loved, I prayed deeply. I wanted to go to morning prayer just to pray for my friends. I li</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t>stened to a lot of sermons too. I thought I was doing all that I could to love and be loved.</w:t>\n      </w:r>\n      <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n      <w:bookmarkEnd w:id=\"0\"/>\n      <w:r>\n        <w:rPr>\n          <w:sz w:val=\"24\"/>\n          <w:szCs w:val=\"24\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> </w:t>\n      </w:r>\n    </w:p>\n</w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n[void]$d.Content.InsertXML($ooxml)\n"}
